$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Source:" citation block (rows 44-51) is being reformatted: the long
# citation sentence that used to live in one hyperlinked cell is split
# across several short lines (with blank spacer lines in between), and the
# trailing long citation is replaced by a second, plain "NIS" line.
#
# Before:
#   44 Source:
#   45 Entreprises privées selon l'activité principale et la tranche de
#      salariés (RNE 2011)
#   46 http://www.ins.nat.tn/...pdf   (hyperlinked)
#   47 p. 69
#   50 NIS
#   51 National Institute of Statistics (NIS), "STATISTIQUES..." (long cite)
#
# After:
#   44 Source:
#   45 (blank)
#   46 Entreprises privées selon l'activité principale et la tranche de
#      salariés (RNE 2011)
#   47 (blank)
#   48 http://www.ins.nat.tn/...pdf   (no longer hyperlinked)
#   49 (blank)
#   50 p. 69
#   53 NIS
#   54 NIS

# Insert three blank rows, pushing rows 45-51 down so each existing line
# gets a blank spacer above it (the new rows inherit the "source" look of
# the row directly above them).
$ws.Rows(45).Insert()
$ws.Rows(47).Insert()
$ws.Rows(49).Insert()

# Row insertion carried the hyperlink along (now anchored at A48) and
# stamped A48/A49 with the HyperLink look (bold/underline/blue); strip the
# link and restore the plain italic "source" citation look used by the
# rest of the block (matching A44/A46/A50).
$ws.Range("A48").Hyperlinks.Delete()
$ws.Range("A48:A49").Font.Underline = $false
$ws.Range("A48:A49").Font.Bold = $false
$ws.Range("A48:A49").Font.Italic = $true
$ws.Range("A48:A49").Font.ColorIndex = -4105

# The long citation (previously at A51, now shifted to A54) is replaced by
# a short, plain repeat of "NIS".
$ws.Range("A54").Value = "NIS"
